# [PHOENIX-5854] refactoring Trade License Module
#
# Applies the changes captured by the target OOXML diff:
#  - tradeOwnerDetails (sheet1): update the sample mobile number, narrow the
#    mobileNumber column, and make this sheet the active/selected one.
#  - tradeLocationDetails (sheet2): insert two new columns (locality, ward)
#    between propertyAssessmentDetails and ownershipType, with sample data,
#    and refresh the sheet selection.
#  - tradeDetails (sheet3): it is no longer the active sheet; update its
#    lingering selection.

$wb = $excel.ActiveWorkbook

$wsOwner    = $wb.Worksheets.Item("tradeOwnerDetails")
$wsLocation = $wb.Worksheets.Item("tradeLocationDetails")
$wsDetails  = $wb.Worksheets.Item("tradeDetails")

# ---------------------------------------------------------------------------
# tradeOwnerDetails (sheet1)
# ---------------------------------------------------------------------------

# Replace the sample mobile number used by the functional test data.
$wsOwner.Cells.Item(2, 3).Value = 9036544535

# The mobileNumber column narrows and becomes a "best fit" custom width.
$wsOwner.Columns.Item(3).ColumnWidth = 12.02

# ---------------------------------------------------------------------------
# tradeLocationDetails (sheet2) - insert "locality" and "ward" columns
# ---------------------------------------------------------------------------

# Shift the ownershipType column (and its data) two columns to the right,
# opening up C:D for the two new fields.
$wsLocation.Range("C1:D2").Insert(-4161)

$wsLocation.Range("C1").Value = "locality"
$wsLocation.Range("D1").Value = "ward"
$wsLocation.Range("C2").Value = "kotha peta"
$wsLocation.Range("D2").Value = "Revenue Ward No 41"

# Matches the refreshed text formatting captured on the dataName cell.
$wsLocation.Range("A2").NumberFormat = "@"

# Approximate the new column widths for the inserted + shifted columns.
$wsLocation.Columns.Item(2).ColumnWidth = 24.28515625
$wsLocation.Columns.Item(3).ColumnWidth = 9.7109375
$wsLocation.Columns.Item(4).ColumnWidth = 18.7109375

# ---------------------------------------------------------------------------
# Selections / active sheet
# ---------------------------------------------------------------------------
# tradeOwnerDetails becomes the active tab (activeTab is no longer forced to
# tradeDetails); update the lingering selections on the other two sheets
# first so that the very last Select() call below is what sticks as active.

$wsLocation.Range("G18").Select()
$wsDetails.Range("G9").Select()

$wsOwner.Range("D10:D11").Select()
